# Delete the row for "HBA" (Hobart, Australia), which causes all following
# rows in the colo list to shift up by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(320).Delete()
